# update new orleans xlsx files
#
# 1. Reorder the worksheet tabs so "review_info" comes before "hotel_info".
# 2. Add a new "State" column to hotel_info (inserted right after
#    Hotel_Name, before City) populated with "Louisiana" for the existing
#    hotel row.

$wb = $excel.ActiveWorkbook

# Move review_info in front of hotel_info so it becomes the first tab.
$reviewSheet = $wb.Worksheets.Item("review_info")
$reviewSheet.Move($wb.Worksheets.Item(1))

# Re-fetch hotel_info by name (NOT the variable captured before the move --
# worksheet references track tab position, so a stale reference would now
# point at review_info instead).
$hotelSheet = $wb.Worksheets.Item("hotel_info")

# Insert a new "State" column in hotel_info between Hotel_Name (col B) and
# City (the old col C, now shifted to col D).
$hotelSheet.Columns.Item(3).Insert()
$hotelSheet.Range("C1").Value = "State"
$hotelSheet.Range("C2").Value = "Louisiana"
